$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05040950907422825
$ws.Range("D2").Value = 0.01094241228024373
$ws.Range("E2").Value = 2.60570920787319
$ws.Range("F2").Value = 0.3535886159855934
$ws.Range("G2").Value = 0.002349316804360842
$ws.Range("M2").Value = 10.73766856433434
$ws.Range("O2").Value = 1.045566560175445
$ws.Range("C3").Value = 0.04465626906139164
$ws.Range("D3").Value = 0.009650969053296876
$ws.Range("E3").Value = 2.268374438788129
$ws.Range("F3").Value = 0.3559578068447848
$ws.Range("G3").Value = 0.002355671204779274
$ws.Range("M3").Value = 9.382611207486832
$ws.Range("O3").Value = 1.071878505284758
$ws.Range("C4").Value = 0.04114654460288136
$ws.Range("D4").Value = 0.008877739012717711
$ws.Range("E4").Value = 2.061642785782851
$ws.Range("F4").Value = 0.3586509693295596
$ws.Range("G4").Value = 0.002359727480142701
$ws.Range("M4").Value = 8.549468863310494
$ws.Range("O4").Value = 1.092498922142084
$ws.Range("C5").Value = 0.03972176177171605
$ws.Range("D5").Value = 0.008567236292073233
$ws.Range("E5").Value = 1.977468670141349
$ws.Range("F5").Value = 0.3600501848854591
$ws.Range("G5").Value = 0.002361419660900487
$ws.Range("M5").Value = 8.209583608508126
$ws.Range("O5").Value = 1.101988250983311
$ws.Range("C6").Value = 0.03948549958549563
$ws.Range("D6").Value = 0.00851594483003737
$ws.Range("E6").Value = 1.963495120641312
$ws.Range("F6").Value = 0.3603004674694574
$ws.Range("G6").Value = 0.002361703024023022
$ws.Range("M6").Value = 8.15312097578601
$ws.Range("O6").Value = 1.103628510007326
$ws.Range("C7").Value = 0.0411273077338592
$ws.Range("D7").Value = 0.008873533325726157
$ws.Range("E7").Value = 2.060507333224308
$ws.Range("F7").Value = 0.3586686306122857
$ws.Range("G7").Value = 0.002359750142138239
$ws.Range("M7").Value = 8.544886672702376
$ws.Range("O7").Value = 1.092622547482819
$ws.Range("C8").Value = 0.04842089385668658
$ws.Range("D8").Value = 0.01049280307873346
$ws.Range("E8").Value = 2.489296668993205
$ws.Range("F8").Value = 0.3541425747761764
$ws.Range("G8").Value = 0.002351475907438469
$ws.Range("M8").Value = 10.2706233464948
$ws.Range("O8").Value = 1.053690140137604
$ws.Range("C9").Value = 0.0629178454185535
$ws.Range("D9").Value = 0.01384235374686682
$ws.Range("E9").Value = 3.334721912355803
$ws.Range("F9").Value = 0.3555394630784576
$ws.Range("G9").Value = 0.002336461405685986
$ws.Range("M9").Value = 13.65040725543832
$ws.Range("O9").Value = 1.014441380901218
$ws.Range("C10").Value = 0.07370754594207085
$ws.Range("D10").Value = 0.01643593458374681
$ws.Range("E10").Value = 3.96092762727892
$ws.Range("F10").Value = 0.3634606177371467
$ws.Range("G10").Value = 0.002326145213537379
$ws.Range("M10").Value = 16.13816614494073
$ws.Range("O10").Value = 1.010592119286486
$ws.Range("C11").Value = 0.07865092151085662
$ws.Range("D11").Value = 0.01765098275109267
$ws.Range("E11").Value = 4.247476084258722
$ws.Range("F11").Value = 0.3687042696817713
$ws.Range("G11").Value = 0.002321601893500613
$ws.Range("M11").Value = 17.27274156481701
$ws.Range("O11").Value = 1.014803359256945
$ws.Range("C12").Value = 0.08052828017510194
$ws.Range("D12").Value = 0.01811672086920879
$ws.Range("E12").Value = 4.356277938613573
$ws.Range("F12").Value = 0.3709380758484144
$ws.Range("G12").Value = 0.002319902498854083
$ws.Range("M12").Value = 17.70295277217264
$ws.Range("O12").Value = 1.017301759571865
$ws.Range("C13").Value = 0.08012371081301239
$ws.Range("D13").Value = 0.01801615624838604
$ws.Range("E13").Value = 4.332831673984003
$ws.Range("F13").Value = 0.3704457485017087
$ws.Range("G13").Value = 0.002320267564574052
$ws.Range("M13").Value = 17.6102709806936
$ws.Range("O13").Value = 1.016722734448308
$ws.Range("C14").Value = 0.07880526237299534
$ws.Range("D14").Value = 0.01768918319609014
$ws.Range("E14").Value = 4.256421120383891
$ws.Range("F14").Value = 0.3688830009195385
$ws.Range("G14").Value = 0.002321461663775753
$ws.Range("M14").Value = 17.3081228789905
$ws.Range("O14").Value = 1.014990520698319
$ws.Range("C15").Value = 0.07799839023530808
$ws.Range("D15").Value = 0.01748965275975678
$ws.Range("E15").Value = 4.209657059587528
$ws.Range("F15").Value = 0.3679584619186897
$ws.Range("G15").Value = 0.002322195813209367
$ws.Range("M15").Value = 17.1231278752611
$ws.Range("O15").Value = 1.014048584474068
$ws.Range("C16").Value = 0.07338522149150606
$ws.Range("D16").Value = 0.01635728407019599
$ws.Range("E16").Value = 3.94223907419331
$ws.Range("F16").Value = 0.3631519781038577
$ws.Range("G16").Value = 0.002326445102469432
$ws.Range("M16").Value = 16.06409011922676
$ws.Range("O16").Value = 1.010440799180458
$ws.Range("C17").Value = 0.07056443207879681
$ws.Range("D17").Value = 0.01567201932365947
$ws.Range("E17").Value = 3.778652349183687
$ws.Range("F17").Value = 0.3606320580864661
$ws.Range("G17").Value = 0.002329089902406649
$ws.Range("M17").Value = 15.41525057961525
$ws.Range("O17").Value = 1.009787007607656
$ws.Range("C18").Value = 0.06894526831625569
$ws.Range("D18").Value = 0.01528115182679102
$ws.Range("E18").Value = 3.684716873193963
$ws.Range("F18").Value = 0.3593362954256918
$ws.Range("G18").Value = 0.002330625227568882
$ws.Range("M18").Value = 15.04231616580824
$ws.Range("O18").Value = 1.009969284509026
$ws.Range("C19").Value = 0.06839759976627136
$ws.Range("D19").Value = 0.01514935744742729
$ws.Range("E19").Value = 3.652937155020879
$ws.Range("F19").Value = 0.3589235919255032
$ws.Range("G19").Value = 0.00233114749804814
$ws.Range("M19").Value = 14.91608721833109
$ws.Range("O19").Value = 1.010125510864384
$ws.Range("C20").Value = 0.07086436790828543
$ws.Range("D20").Value = 0.01574462365408635
$ws.Range("E20").Value = 3.796050016567449
$ws.Range("F20").Value = 0.3608843204698928
$ws.Range("G20").Value = 0.002328806902309966
$ws.Range("M20").Value = 15.48429255328728
$ws.Range("O20").Value = 1.009798491126219
$ws.Range("C21").Value = 0.07919237283314828
$ws.Range("D21").Value = 0.01778506574900263
$ws.Range("E21").Value = 4.278856401999349
$ws.Range("F21").Value = 0.3693351828021889
$ws.Range("G21").Value = 0.002321110359947632
$ws.Range("M21").Value = 17.39685424505313
$ws.Range("O21").Value = 1.015474409292068
$ws.Range("C22").Value = 0.08466691874335197
$ws.Range("D22").Value = 0.01915169384837867
$ws.Range("E22").Value = 4.59612777059408
$ws.Range("F22").Value = 0.3763107667697625
$ws.Range("G22").Value = 0.002316202766369164
$ws.Range("M22").Value = 18.65024221764463
$ws.Range("O22").Value = 1.024474255923622
$ws.Range("C23").Value = 0.08174202380692464
$ws.Range("D23").Value = 0.01841907551798982
$ws.Range("E23").Value = 4.426618119971693
$ws.Range("F23").Value = 0.3724506806521646
$ws.Range("G23").Value = 0.002318810982607245
$ws.Range("M23").Value = 17.98091679736842
$ws.Range("O23").Value = 1.019171000089074
$ws.Range("C24").Value = 0.07072875901737063
$ws.Range("D24").Value = 0.01571178962142739
$ws.Range("E24").Value = 3.78818418379592
$ws.Range("F24").Value = 0.3607697970931554
$ws.Range("G24").Value = 0.002328934800503592
$ws.Range("M24").Value = 15.45307840193578
$ws.Range("O24").Value = 1.009791564396636
$ws.Range("C25").Value = 0.05897318681807917
$ws.Range("D25").Value = 0.012914943153433
$ws.Range("E25").Value = 3.10529981578776
$ws.Range("F25").Value = 0.3539977218348227
$ws.Range("G25").Value = 0.002340395904093154
$ws.Range("M25").Value = 12.73584327026504
$ws.Range("O25").Value = 1.020847501603299

Write-Host "Updated 380 kV case values"
